$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TP")

$ws.Range("F7").Value = 324.15
$ws.Range("G7").Value = 330.7
$ws.Range("H7").Value = 320.05
$ws.Range("I7").Value = 328.55
$ws.Range("J7").Value = 320.15

$ws.Range("G9").Value = 325.25
$ws.Range("H9").Value = 312.7
$ws.Range("I9").Value = 322.9

$ws.Range("G10").Value = 328.55
$ws.Range("H10").Value = 322.6
$ws.Range("I10").Value = 327.65

$ws.Range("G11").Value = 330
$ws.Range("H11").Value = 326.8
$ws.Range("I11").Value = 327.1

$ws.Range("G12").Value = 329
$ws.Range("H12").Value = 325
$ws.Range("I12").Value = 327.85

$ws.Range("G13").Value = 328.4
$ws.Range("H13").Value = 325.55
$ws.Range("I13").Value = 325.8

$ws.Range("G14").Value = 328.15
$ws.Range("H14").Value = 325.75
$ws.Range("I14").Value = 327.75

$ws.Range("G15").Value = 328.9
$ws.Range("H15").Value = 327.1
$ws.Range("I15").Value = 328.2

$ws.Range("G16").Value = 330.7
$ws.Range("H16").Value = 328.05
$ws.Range("I16").Value = 329.1

$ws.Range("G17").Value = 329.7
$ws.Range("H17").Value = 327
$ws.Range("I17").Value = 327.95

$ws.Range("G18").Value = 328.9
$ws.Range("H18").Value = 325
$ws.Range("I18").Value = 325.65

$ws.Range("G19").Value = 327.45
$ws.Range("H19").Value = 324
$ws.Range("I19").Value = 327.25

$ws.Range("G20").Value = 328.35
$ws.Range("H20").Value = 326.8
$ws.Range("I20").Value = 328.1

$ws.Range("G21").Value = 329.2
$ws.Range("H21").Value = 327.95
$ws.Range("I21").Value = 328.95
